# daily auto push: 2026-02-26 03:09 UTC
#
# The data sheet is an append-only log of (date, weekday, hour, value)
# samples. A sample for 2026/02/26 at hour 9 (value 182) was missing from
# between the existing 2026/02/26 06:00 row (886) and the 2026/12/29 13:00
# row (then-887). Insert it as a new row 887, pushing rows 887-928 down to
# 888-929 (dimension grows from D928 to D929).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 887; Excel shifts 887..928 -> 888..929 and the
# sheet's used-range dimension grows automatically.
$ws.Rows.Item(887).Insert()

# Column A holds the date as plain text ("2026/12/29", etc.), not a real
# date value. Pre-format the cell as Text so Excel doesn't auto-convert the
# slash-separated string into a date serial, then drop back to the sheet's
# normal (unstyled) formatting so the new row matches its neighbours.
$ws.Range("A887").NumberFormat = "@"
$ws.Range("A887").Value = "2026/02/26"
$ws.Range("A887").ClearFormats()

$ws.Range("B887").Value = "木"
$ws.Range("C887").Value = 9
$ws.Range("D887").Value = 182
